$d = $word.ActiveDocument

# A register ("reg1") is allowed to have no description in IP-XACT.
# The test fixture is updated so that reg1 has no description:
#  - the table-cell paragraph that held the placeholder description text
#    becomes an empty paragraph (cell must keep at least one paragraph)
#  - the standalone "Definition" paragraph that held the placeholder
#    description text (underneath the "Description" term) is removed
#    entirely

$target = "write something useful for reg1"

foreach ($p in @($d.Paragraphs)) {
    $txt = $p.Range.Text
    if ($null -eq $txt) { continue }

    # Range.Text carries a trailing paragraph mark (CR) and, for the
    # last/only paragraph in a table cell, a trailing cell mark (BEL) too.
    $clean = $txt.TrimEnd([char]13, [char]7)
    if ($clean -ne $target) { continue }

    $styleName = $p.Style.NameLocal

    if ($styleName -eq "Compact") {
        # Table-cell paragraph -> blank it out, keeping an empty <w:p/>
        # (delete the run text but not the paragraph mark, then drop the
        # paragraph-level formatting so no <w:pPr> survives)
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Delete()
        $p.Style = "Normal"
        $p.Alignment = 0
    }
    elseif ($styleName -eq "Definition") {
        # Standalone definition paragraph -> remove it completely
        $p.Range.Delete()
    }
}
